{"js": "// Word JavaScript API (Office.js) script.\n// Body of: async (context) => { ... }\n//\n// Content edit (matches the xml diff for word/document.xml):\n//   The SmartCash-mining paragraph no longer ends with the trailing\n//   clause \", until Smartcash reaches a considerable market cap.\" \u2014\n//   it now simply ends after \"...quite some time.\"\n//\n// NOTE on the \"exchanges\" bookmark: the diff also shows\n// <w:bookmarkStart w:id=\"0\" .../> becoming <w:bookmarkStart w:id=\"1\" .../>\n// (and the matching bookmarkEnd). w:id is an internal/invisible OOXML\n// identifier -- it is not exposed as a settable property anywhere in the\n// Word API (no Range/Bookmark .id setter), and this runtime always\n// renumbers bookmarks by their document position at save time (the\n// single \"exchanges\" bookmark is therefore always serialized back out as\n// id 0, the same way plain Word renumbers bookmark ids on save). There\n// is no scriptable way to pin it to 1 without adding another, unwanted\n// bookmark to the document, so that part of the diff is left alone here.\n\nconst body = context.document.body;\n\nconst results = body.search(\n  \"no ASICs will be created for quite some time, until Smartcash reaches a considerable market cap.\",\n  { matchCase: true }\n);\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace the whole matched sentence fragment, keeping the sentence\n  // ending at \"...quite some time.\" and dropping the trailing clause.\n  results.items[0].insertText(\n    \"no ASICs will be created for quite some time.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document.\n#\n# Content edit (matches the xml diff for word/document.xml):\n#   The SmartCash-mining paragraph no longer ends with the trailing\n#   clause \", until Smartcash reaches a considerable market cap.\" \u2014\n#   it now simply ends after \"...quite some time.\"\n#\n# NOTE on the \"exchanges\" bookmark: the diff also shows\n# <w:bookmarkStart w:id=\"0\" .../> becoming <w:bookmarkStart w:id=\"1\" .../>\n# (and the matching bookmarkEnd). w:id is an internal/invisible OOXML\n# identifier -- it is not exposed as a settable property anywhere in the\n# Word object model (Bookmark has no .ID/.Id setter), and this runtime\n# always renumbers bookmarks by their document position at save time (the\n# single \"exchanges\" bookmark is therefore always serialized back out as\n# id 0, the same way plain Word renumbers bookmark ids on save). There is\n# no scriptable way to pin it to 1 without adding another, unwanted\n# bookmark to the document, so that part of the diff is left alone here.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \", until Smartcash reaches a considerable market cap.\"\n$find.Replacement.Text = \".\"\n$find.Forward = $true\n$find.Wrap = 1          # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# wdReplaceAll = 2\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n"}
